$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Copy the formatting of the row above (row 12) for the label column (A)
# so the new year label picks up the same bold/centered/bordered style (s="1").
$ws.Cells.Item(12, 1).Copy($ws.Cells.Item($row, 1))

# Column A: year label
$ws.Cells.Item($row, 1).Value = "2021年"

# Numeric columns
$ws.Cells.Item($row, 2).Value = 5942
$ws.Cells.Item($row, 3).Value = 7296
$ws.Cells.Item($row, 4).Value = 465760

# Column E is blank (an explicit empty text cell, matching the other rows
# above it). Writing a leading apostrophe forces a text cell with empty
# content, then clearing the style drops the "quote prefix" formatting flag
# that the apostrophe trick leaves behind so the cell ends up unstyled, just
# like the existing blank cells in this column.
$ws.Cells.Item($row, 5).Value = "'"
$ws.Cells.Item($row, 5).Style = "Normal"

$ws.Cells.Item($row, 6).Value = 975714
$ws.Cells.Item($row, 7).Value = 2168130
$ws.Cells.Item($row, 8).Value = 762147
$ws.Cells.Item($row, 9).Value = 227945
$ws.Cells.Item($row, 10).Value = 448219
$ws.Cells.Item($row, 11).Value = 920611
$ws.Cells.Item($row, 12).Value = 1613126

# Column M is blank as well, same trick as column E.
$ws.Cells.Item($row, 13).Value = "'"
$ws.Cells.Item($row, 13).Style = "Normal"

$ws.Cells.Item($row, 14).Value = 21050
$ws.Cells.Item($row, 15).Value = 131720
$ws.Cells.Item($row, 16).Value = 555702
$ws.Cells.Item($row, 17).Value = 2533737
$ws.Cells.Item($row, 18).Value = 10764
$ws.Cells.Item($row, 19).Value = 59996
